# Updated preprocessing of patientValues (vitals) for extraction of 08 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the units for the newly-documented vitals (FiO2, heart frequency, respiratory frequency)
$ws.Range("C5").Value = "%"
$ws.Range("C8").Value = "bpm"
$ws.Range("C16").Value = "/min"

# Reflect the reviewer's final scroll position / zoom / selection on the sheet
$win = $excel.ActiveWindow
[void]$ws.Range("B17").Select()
$win.Zoom = 169
[void]$ws.Range("C27").Select()
